# Applies the commit's edits to the workbook:
#  - Metadata sheet: update URL, Version, Date, Publisher values
#  - Elements sheet: clear the Constraint(s) value for the root "Extension" row

$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/dental-fully-insured-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; column AI is "Constraint(s)" -> clear it
$elements.Range("AI2").Value = ""
